# Auto-generated edit script: applies weekly permutation of rows 2-25
# and appends a new row 26 to the Arandano (blue) / Macroferia Talca sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colIndex = @{"D"=4;"L"=12;"M"=13;"N"=14;"O"=15;"P"=16;"Q"=17;"R"=18;"S"=19;"T"=20}

# 1) Snapshot the current (pre-edit) values of the variable columns for every data row (2-25).
$snapshot = @{}
for ($r = 2; $r -le 25; $r++) {
    $rowData = @{}
    foreach ($col in $colIndex.Keys) {
        $rowData[$col] = $ws.Cells.Item($r, $colIndex[$col]).Value2
    }
    $snapshot[$r] = $rowData
}

# 2) Row mapping: new row number -> source row number (values copied from snapshot).
$rowMap = @{
    2 = 7
    3 = 21
    4 = 4
    5 = 15
    6 = 24
    7 = 25
    8 = 9
    9 = 6
    10 = 12
    11 = 13
    12 = 14
    13 = 5
    14 = 8
    15 = 22
    16 = 23
    17 = 18
    18 = 19
    19 = 3
    20 = 16
    21 = 17
    22 = 11
    23 = 10
    24 = 20
    25 = 2
}

# 3) Rewrite rows 2-25 using the snapshotted source-row values.
foreach ($newRow in $rowMap.Keys) {
    $srcRow = $rowMap[$newRow]
    $src = $snapshot[$srcRow]
    foreach ($col in $colIndex.Keys) {
        $ws.Cells.Item($newRow, $colIndex[$col]).Value = $src[$col]
    }
}

# 4) Append the new row 26 (fresh weekly data point).
$ws.Cells.Item(26, $colIndex["D"]).Value = 44517
$ws.Cells.Item(26, $colIndex["L"]).Value = "Primera"
$ws.Cells.Item(26, $colIndex["M"]).Value = 20
$ws.Cells.Item(26, $colIndex["N"]).Value = 5000
$ws.Cells.Item(26, $colIndex["O"]).Value = 5000
$ws.Cells.Item(26, $colIndex["P"]).Value = 5000
$ws.Cells.Item(26, $colIndex["Q"]).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(26, $colIndex["R"]).Value = "Provincia de Linares"
$ws.Cells.Item(26, $colIndex["S"]).Value = 2500
$ws.Cells.Item(26, $colIndex["T"]).Value = 2

# Row 26 shares the same fixed (non-varying) descriptor columns as every other row.
$ws.Cells.Item(26, 1).Value = $ws.Cells.Item(25, 1).Value2
$ws.Cells.Item(26, 2).Value = $ws.Cells.Item(25, 2).Value2
$ws.Cells.Item(26, 3).Value = $ws.Cells.Item(25, 3).Value2
$ws.Cells.Item(26, 5).Value = $ws.Cells.Item(25, 5).Value2
$ws.Cells.Item(26, 6).Value = $ws.Cells.Item(25, 6).Value2
$ws.Cells.Item(26, 7).Value = $ws.Cells.Item(25, 7).Value2
$ws.Cells.Item(26, 8).Value = $ws.Cells.Item(25, 8).Value2
$ws.Cells.Item(26, 9).Value = $ws.Cells.Item(25, 9).Value2
$ws.Cells.Item(26, 10).Value = $ws.Cells.Item(25, 10).Value2
$ws.Cells.Item(26, 11).Value = $ws.Cells.Item(25, 11).Value2

# 5) Match the date-formatted style used by every other cell in column D.
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(25, 4).NumberFormat

